$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" column header (H1) - copy the formatting of the existing
# header cells (G1 "sum") so it keeps the same bold/border/center style,
# then set its text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Save values per row (2..22)
$saveValues = @{
    2  = 1
    3  = 0
    4  = 0
    5  = 0
    6  = 0
    7  = 1
    8  = 1
    9  = 0
    10 = 1
    11 = 0
    12 = 0
    13 = 0
    14 = 1
    15 = 0
    16 = 0
    17 = 0
    18 = 0
    19 = 0
    20 = 1
    21 = 0
    22 = 0
}

foreach ($row in $saveValues.Keys) {
    $ws.Range("H$row").Value = $saveValues[$row]
}
